# Staging.Institution.xlsx: the four header-row labels in row 2 were
# rotated - the column that used to read "Institution_ID" now reads
# "BusinessKey", "Name" now reads "Institution_ID", and the old
# "BusinessKey" column now reads "Name" ("Code" in column B is unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "BusinessKey"
$ws.Range("B2").Value = "Code"
$ws.Range("C2").Value = "Institution_ID"
$ws.Range("D2").Value = "Name"
